$d = $word.ActiveDocument

# "... 6 Level designs..." -> "... 5 Level designs..."
# (Keep the match inside the un-bolded run that already holds "6 Level
# designs...", so the replacement text inherits that same formatting
# instead of merging into the neighbouring bold "Sam McMillan:" run.)
$d.Content.Find.Execute(
    "6 Level designs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5 Level designs", 2) | Out-Null

# Move the "_GoBack" bookmark from the end of the "Feedback Received:"
# paragraph to sit right after the "5" (between "5" and " Level designs...")
# in the "Sam McMillan:" paragraph. Adding a bookmark with the same name
# repositions it (a bookmark name is unique per document).
$r = $d.Content
$r.Find.Execute(
    "Sam McMillan: 5", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null

$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
